# Auto-generated edit script: updates profit/price figures across multiple
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 71.59999999999999
$ws.Range("I4").Value = 71.59999999999999
$ws.Range("K4").Value = 71.59999999999999
$ws.Range("M4").Value = 42.40000000000001
$ws.Range("H17").Value = 5925.8423
$ws.Range("J17").Value = 6523
$ws.Range("L17").Value = 19569
$ws.Range("N17").Value = -19905
$ws.Range("H62").Value = 6713.5713
$ws.Range("I62").Value = 4664.8335
$ws.Range("K62").Value = 4664.8335
$ws.Range("M62").Value = -4040.8335
$ws.Range("H65").Value = 6713.5713
$ws.Range("I65").Value = 4664.8335
$ws.Range("K65").Value = 23324.1675
$ws.Range("M65").Value = -20204.1675
$ws.Range("H96").Value = 1726.0769
$ws.Range("I96").Value = 1125.6364
$ws.Range("K96").Value = 3376.9092
$ws.Range("M96").Value = -2003.9092
$ws.Range("H103").Value = 2597
$ws.Range("J103").Value = 2597
$ws.Range("L103").Value = 7791
$ws.Range("N103").Value = -8963
$ws.Range("H107").Value = 2454.5454
$ws.Range("I107").Value = 2333.889
$ws.Range("K107").Value = 2333.889
$ws.Range("M107").Value = -413.8890000000001
$ws.Range("H127").Value = 993.13336
$ws.Range("I127").Value = 645.8
$ws.Range("J127").Value = 1687.8
$ws.Range("K127").Value = 1937.4
$ws.Range("L127").Value = 5063.4
$ws.Range("M127").Value = 3022.6
$ws.Range("N127").Value = -14983.4
$ws.Range("H132").Value = 5200.5
$ws.Range("I132").Value = 5240.6
$ws.Range("K132").Value = 15721.8
$ws.Range("M132").Value = -13191.8
$ws.Range("H137").Value = 3511.923
$ws.Range("I137").Value = 2973.1667
$ws.Range("K137").Value = 8919.500100000001
$ws.Range("M137").Value = -6369.500100000001
$ws.Range("H138").Value = 3415.6365
$ws.Range("J138").Value = 3558.5862
$ws.Range("L138").Value = 10675.7586
$ws.Range("N138").Value = -20955.7586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 388
$ws.Range("I5").Value = 183
$ws.Range("K5").Value = 183
$ws.Range("M5").Value = -71
$ws.Range("H61").Value = 3466.1428
$ws.Range("I61").Value = 2670.0588
$ws.Range("K61").Value = 2670.0588
$ws.Range("M61").Value = -2458.0588
$ws.Range("H74").Value = 14494221
$ws.Range("I74").Value = 20834678
$ws.Range("K74").Value = 20834678
$ws.Range("M74").Value = -20833804
$ws.Range("H77").Value = 14494221
$ws.Range("I77").Value = 20834678
$ws.Range("K77").Value = 104173390
$ws.Range("M77").Value = -104169022
$ws.Range("H110").Value = 4317
$ws.Range("I110").Value = 405.5
$ws.Range("K110").Value = 405.5
$ws.Range("M110").Value = 1639.5
$ws.Range("H136").Value = 3466.1428
$ws.Range("I136").Value = 2670.0588
$ws.Range("K136").Value = 8010.176399999999
$ws.Range("M136").Value = -5460.176399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 388
$ws.Range("I4").Value = 183
$ws.Range("K4").Value = 183
$ws.Range("M4").Value = -68
$ws.Range("H134").Value = 4593.8887
$ws.Range("I134").Value = 2551.077
$ws.Range("K134").Value = 7653.231000000001
$ws.Range("M134").Value = -5118.231000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H4").Value = 1754500
$ws.Range("J4").Value = 1754500
$ws.Range("L4").Value = 1754500
$ws.Range("N4").Value = -1754724
$ws.Range("H6").Value = 3998.5
$ws.Range("I6").Value = 3998.5
$ws.Range("K6").Value = 3998.5
$ws.Range("M6").Value = -3885.5
$ws.Range("H7").Value = 180.91667
$ws.Range("I7").Value = 74.57143000000001
$ws.Range("J7").Value = 329.8
$ws.Range("K7").Value = 74.57143000000001
$ws.Range("L7").Value = 329.8
$ws.Range("M7").Value = 38.42856999999999
$ws.Range("N7").Value = -555.8
$ws.Range("H16").Value = 1277.55
$ws.Range("I16").Value = 1294.9375
$ws.Range("K16").Value = 1294.9375
$ws.Range("M16").Value = -1007.9375
$ws.Range("H31").Value = 37551.75
$ws.Range("I31").Value = 3714.1765
$ws.Range("J31").Value = 75901
$ws.Range("K31").Value = 3714.1765
$ws.Range("L31").Value = 75901
$ws.Range("M31").Value = -3419.1765
$ws.Range("N31").Value = -76491
$ws.Range("H34").Value = 37551.75
$ws.Range("I34").Value = 3714.1765
$ws.Range("J34").Value = 75901
$ws.Range("K34").Value = 3714.1765
$ws.Range("L34").Value = 75901
$ws.Range("M34").Value = -3512.1765
$ws.Range("N34").Value = -76305
$ws.Range("H55").Value = 15457.5
$ws.Range("I55").Value = 11373
$ws.Range("J55").Value = 17499.75
$ws.Range("K55").Value = 11373
$ws.Range("L55").Value = 17499.75
$ws.Range("M55").Value = -11058
$ws.Range("N55").Value = -18129.75
$ws.Range("H62").Value = 29006
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 29006
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 29006
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -30254
$ws.Range("H65").Value = 29006
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 29006
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 145030
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -151270
$ws.Range("H113").Value = 1277.55
$ws.Range("I113").Value = 1294.9375
$ws.Range("K113").Value = 1294.9375
$ws.Range("M113").Value = 875.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2852
$ws.Range("J17").Value = 6666
$ws.Range("L17").Value = 19998
$ws.Range("N17").Value = -20336
$ws.Range("H109").Value = 3590.8
$ws.Range("I109").Value = 3323.111
$ws.Range("J109").Value = 6000
$ws.Range("K109").Value = 9969.332999999999
$ws.Range("L109").Value = 18000
$ws.Range("M109").Value = -8929.332999999999
$ws.Range("N109").Value = -20080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2049.25
$ws.Range("I102").Value = 1450.9131
$ws.Range("K102").Value = 1450.9131
$ws.Range("M102").Value = 171.0869

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4429.9585
$ws.Range("J132").Value = 5498.154
$ws.Range("L132").Value = 16494.462
$ws.Range("N132").Value = -21554.462
$ws.Range("H136").Value = 3878.625
$ws.Range("I136").Value = 1993.4445
$ws.Range("K136").Value = 5980.333500000001
$ws.Range("M136").Value = -3430.333500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10226
$ws.Range("H57").Value = 52000
$ws.Range("I57").Value = 52000
$ws.Range("K57").Value = 52000
$ws.Range("M57").Value = -51246
$ws.Range("H81").Value = 1990.4546
$ws.Range("I81").Value = 1990.4546
$ws.Range("K81").Value = 3980.9092
$ws.Range("M81").Value = -2919.9092
$ws.Range("H84").Value = 1990.4546
$ws.Range("I84").Value = 1990.4546
$ws.Range("K84").Value = 19904.546
$ws.Range("M84").Value = -14600.546
$ws.Range("H132").Value = 4684.5713
$ws.Range("I132").Value = 4581.087
$ws.Range("J132").Value = 5160.6
$ws.Range("K132").Value = 13743.261
$ws.Range("L132").Value = 15481.8
$ws.Range("M132").Value = -11213.261
$ws.Range("N132").Value = -20541.8

Write-Host "Updated cells: 184 set, 3 cleared across 8 sheets"
